$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values are stored as text, matching the
# original inline-string cell content (many look like plain numbers,
# e.g. "1.002", and would otherwise be auto-converted to numeric values).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.599.39'
$ws.Range('D3').Value = '1.851.80'
$ws.Range('D5').Value = '336.05'
$ws.Range('D6').Value = '1.002'
$ws.Range('D7').Value = '0.4660'
$ws.Range('D8').Value = '0.3907'
$ws.Range('D10').Value = '0.07909'
$ws.Range('D12').Value = '22.33'
$ws.Range('D13').Value = '1.876.83'
$ws.Range('D14').Value = '5.851'
$ws.Range('D15').Value = '7.015'
$ws.Range('D16').Value = '0.06902'
$ws.Range('D17').Value = '87.63'
$ws.Range('D18').Value = '1.001'
$ws.Range('D19').Value = '0.00001004'
$ws.Range('D20').Value = '17.13'
$ws.Range('D21').Value = '1.002'
$ws.Range('D22').Value = '28.598.65'
$ws.Range('D23').Value = '5.401'
$ws.Range('D24').Value = '11.32'
$ws.Range('D25').Value = '2.129'
$ws.Range('D26').Value = '2.067.04'
$ws.Range('D27').Value = '153.44'
$ws.Range('D28').Value = '19.50'
$ws.Range('D29').Value = '6.015'
$ws.Range('D30').Value = '2.027'
$ws.Range('D31').Value = '117.86'
$ws.Range('D32').Value = '0.9740'
$ws.Range('D33').Value = '0.09393'
$ws.Range('D35').Value = '3.480'
$ws.Range('D36').Value = '1.350'
$ws.Range('D37').Value = '0.06167'
$ws.Range('D38').Value = '0.02199'
$ws.Range('D39').Value = '1.159'
$ws.Range('D40').Value = '0.5727'
$ws.Range('D41').Value = '7.644'
$ws.Range('D43').Value = '0.1802'
$ws.Range('D45').Value = '1.248'
$ws.Range('D46').Value = '0.5394'
$ws.Range('D47').Value = '11.75'
$ws.Range('D48').Value = '0.07143'
$ws.Range('D50').Value = '115.78'
$ws.Range('D51').Value = '43.16'

$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('E5').Value = '  +2.98%  '
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('E8').Value = '  -3.84%  '
$ws.Range('E9').Value = '  -3.44%  '
$ws.Range('E10').Value = '  -3.81%  '
$ws.Range('E11').Value = '  -3.06%  '
$ws.Range('E12').Value = '  -6.56%  '
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('E15').Value = '  -4.19%  '
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('E17').Value = '  -4.46%  '
$ws.Range('E19').Value = '  -3.42%  '
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('E23').Value = '  -4.97%  '
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('E28').Value = '  -2.91%  '
$ws.Range('E29').Value = '  -6.60%  '
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('E34').Value = '  -4.40%  '
$ws.Range('E35').Value = '  -2.27%  '
$ws.Range('E36').Value = '  -2.59%  '
$ws.Range('E37').Value = '  -3.55%  '
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('E40').Value = '  -4.11%  '
$ws.Range('E41').Value = '  -3.32%  '
$ws.Range('E42').Value = '  -5.30%  '
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('E44').Value = '  -5.84%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('E47').Value = '  -5.24%  '
$ws.Range('E48').Value = '  -5.31%  '
$ws.Range('E50').Value = '  -3.01%  '
$ws.Range('E51').Value = '  +1.46%  '
